$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2: swap "mukes" info for "shankar" info
$ws.Range("A2").Value = "shankar"
$ws.Range("C2").Value = "shankar1217"

# Remove the old hyperlink on B2 (mukeshse1@gmail.com) before re-adding the new one,
# then restore the Hyperlink cell style that Delete()/Add() disturbs.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:shankar61@gmail.com", [Type]::Missing, [Type]::Missing, "shankar61@gmail.com")
$ws.Range("B2").Style = "Hyperlink"

# Clear the Skills value that used to sit in F2 - it now belongs to the new row
$ws.Range("F2").ClearContents()

# Add the new second entrant in row 3
$ws.Range("A3").Value = "prem"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:prem15@gmail.com", [Type]::Missing, [Type]::Missing, "prem15@gmail.com")
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("C3").Value = "prem1217"
$ws.Range("D3").Value = 12345456
$ws.Range("E3").Value = 12345456
$ws.Range("F3").Value = "CSS,Java,python"

# Match the saved cursor position recorded in the workbook
[void]$ws.Range("C10").Select()
